$wb = $excel.ActiveWorkbook

# --- packages: bump version/date in the description ---
$packages = $wb.Worksheets.Item("packages")
$packages.Range("C2").Value = "Mapping tables for processing raw data into unified model terminology (v2.1.0, 2022-06-29)"

# --- attributes: tags column (G) now shows the bare NCIT identifier, no more hyperlink ---
$attributes = $wb.Worksheets.Item("attributes")
$attributes.Hyperlinks.Delete()
$attributes.Range("G2").Value = "NCIT_C25516"
$attributes.Range("G3").Value = "NCIT_C65107"
$attributes.Range("G4").Value = "NCIT_C25415"
$attributes.Range("G2:G4").Style = "Normal"

# --- tags: rebuilt model ---
# identifier/label now use the bare NCIT_xxxxx id (identifier loses its hyperlink & style,
# label switches from the "NCIT:xxxxx" form to the "NCIT_xxxxx" form); objectIRI keeps the
# full purl.obolibrary.org IRI + its hyperlink.
$tags = $wb.Worksheets.Item("tags")
$tags.Hyperlinks.Delete()

$tags.Range("A2").Value = "NCIT_C25415"
$tags.Range("B2").Value = "NCIT_C25415"
$tags.Range("C2").Value = "http://purl.obolibrary.org/obo/NCIT_C25415"

$tags.Range("A3").Value = "NCIT_C25516"
$tags.Range("B3").Value = "NCIT_C25516"
$tags.Range("C3").Value = "http://purl.obolibrary.org/obo/NCIT_C25516"

$tags.Range("A4").Value = "NCIT_C65107"
$tags.Range("B4").Value = "NCIT_C65107"
$tags.Range("C4").Value = "http://purl.obolibrary.org/obo/NCIT_C65107"

$tags.Range("A2:A4").Style = "Normal"

# re-create the surviving hyperlinks (objectIRI + relationIRI columns) in ref order
$tags.Hyperlinks.Add($tags.Range("C2"), "http://purl.obolibrary.org/obo/NCIT_C25415") | Out-Null
$tags.Hyperlinks.Add($tags.Range("F2"), "http://molgenis.org", "isAssociatedWith") | Out-Null
$tags.Hyperlinks.Add($tags.Range("C3"), "http://purl.obolibrary.org/obo/NCIT_C25516") | Out-Null
$tags.Hyperlinks.Add($tags.Range("F3"), "http://molgenis.org", "isAssociatedWith") | Out-Null
$tags.Hyperlinks.Add($tags.Range("C4"), "http://purl.obolibrary.org/obo/NCIT_C65107") | Out-Null
$tags.Hyperlinks.Add($tags.Range("F4"), "http://molgenis.org", "isAssociatedWith") | Out-Null

$tags.Range("C2").Style = "Hyperlink"
$tags.Range("C3").Style = "Hyperlink"
$tags.Range("C4").Style = "Hyperlink"
$tags.Range("F2").Style = "Hyperlink"
$tags.Range("F3").Style = "Hyperlink"
$tags.Range("F4").Style = "Hyperlink"
